$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7464590
$ws.Range("J17").Value = 7464590
$ws.Range("L17").Value = 22393770
$ws.Range("N17").Value = -22394106
$ws.Range("H112").Value = 1888.975
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1911.7693
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 5735.3079
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -7951.3079
$ws.Range("H129").Value = 1023.7553
$ws.Range("I129").Value = 2606
$ws.Range("J129").Value = 934.8652
$ws.Range("K129").Value = 7818
$ws.Range("L129").Value = 2804.5956
$ws.Range("M129").Value = -2818
$ws.Range("N129").Value = -12804.5956
$ws.Range("H132").Value = 66228.2
$ws.Range("I132").Value = 9768
$ws.Range("J132").Value = 179148.6
$ws.Range("K132").Value = 29304
$ws.Range("L132").Value = 537445.8
$ws.Range("M132").Value = -26774
$ws.Range("N132").Value = -542505.8
$ws.Range("H137").Value = 2082967
$ws.Range("I137").Value = 5918096.5
$ws.Range("J137").Value = 5605.1665
$ws.Range("K137").Value = 17754289.5
$ws.Range("L137").Value = 16815.4995
$ws.Range("M137").Value = -17751739.5
$ws.Range("N137").Value = -21915.4995
$ws.Range("H138").Value = 1534.38
$ws.Range("I138").Value = 956.26666
$ws.Range("J138").Value = 1782.1428
$ws.Range("K138").Value = 2868.79998
$ws.Range("L138").Value = 5346.428400000001
$ws.Range("M138").Value = 2271.20002
$ws.Range("N138").Value = -15626.4284

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 44916.668
$ws.Range("J37").Value = 44916.668
$ws.Range("L37").Value = 44916.668
$ws.Range("N37").Value = -45462.668
$ws.Range("H61").Value = 2724.8462
$ws.Range("I61").Value = 1741
$ws.Range("K61").Value = 1741
$ws.Range("M61").Value = -1529
$ws.Range("H63").Value = 2880.9375
$ws.Range("I63").Value = 2059.5833
$ws.Range("K63").Value = 2059.5833
$ws.Range("M63").Value = -1373.5833
$ws.Range("H66").Value = 2880.9375
$ws.Range("I66").Value = 2059.5833
$ws.Range("K66").Value = 10297.9165
$ws.Range("M66").Value = -6865.916499999999
$ws.Range("H74").Value = 2067.8572
$ws.Range("I74").Value = 1778.6
$ws.Range("J74").Value = 3514.1428
$ws.Range("K74").Value = 1778.6
$ws.Range("L74").Value = 3514.1428
$ws.Range("M74").Value = -904.5999999999999
$ws.Range("N74").Value = -5262.1428
$ws.Range("H77").Value = 2067.8572
$ws.Range("I77").Value = 1778.6
$ws.Range("J77").Value = 3514.1428
$ws.Range("K77").Value = 8893
$ws.Range("L77").Value = 17570.714
$ws.Range("M77").Value = -4525
$ws.Range("N77").Value = -26306.714
$ws.Range("H132").Value = 20836128
$ws.Range("I132").Value = 27779836
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 83339508
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -83336978
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 2724.8462
$ws.Range("I136").Value = 1741
$ws.Range("K136").Value = 5223
$ws.Range("M136").Value = -2673

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H35").Value = 35423.332
$ws.Range("J35").Value = 35423.332
$ws.Range("L35").Value = 35423.332
$ws.Range("N35").Value = -36043.332
$ws.Range("H82").Value = 11073.571
$ws.Range("H85").Value = 11073.571
$ws.Range("H105").Value = 2865.4092
$ws.Range("I105").Value = 1511.8
$ws.Range("J105").Value = 3993.4167
$ws.Range("K105").Value = 1511.8
$ws.Range("L105").Value = 3993.4167
$ws.Range("M105").Value = 235.2
$ws.Range("N105").Value = -7487.4167
$ws.Range("H134").Value = 2729.7612
$ws.Range("I134").Value = 2004.4857
$ws.Range("K134").Value = 6013.4571
$ws.Range("M134").Value = -3478.4571

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 333380000
$ws.Range("J23").Value = 70010
$ws.Range("L23").Value = 70010
$ws.Range("N23").Value = -70490
$ws.Range("H27").Value = 333380000
$ws.Range("J27").Value = 70010
$ws.Range("L27").Value = 70010
$ws.Range("N27").Value = -70394
$ws.Range("H58").Value = 3303.6875
$ws.Range("I58").Value = 1965
$ws.Range("J58").Value = 5534.8335
$ws.Range("K58").Value = 1965
$ws.Range("L58").Value = 5534.8335
$ws.Range("M58").Value = -1762
$ws.Range("N58").Value = -5940.8335
$ws.Range("H122").Value = 61540.7
$ws.Range("I122").Value = 87296.42999999999
$ws.Range("J122").Value = 1444
$ws.Range("K122").Value = 261889.29
$ws.Range("L122").Value = 4332
$ws.Range("M122").Value = -259439.29
$ws.Range("N122").Value = -9232
$ws.Range("H132").Value = 64218.566
$ws.Range("I132").Value = 1756.2727
$ws.Range("J132").Value = 121475.664
$ws.Range("K132").Value = 5268.8181
$ws.Range("L132").Value = 364426.992
$ws.Range("M132").Value = -2738.8181
$ws.Range("N132").Value = -369486.992
$ws.Range("H134").Value = 1108215.2
$ws.Range("I134").Value = 1566768.5
$ws.Range("J134").Value = 282819.2
$ws.Range("K134").Value = 4700305.5
$ws.Range("L134").Value = 848457.6000000001
$ws.Range("M134").Value = -4697770.5
$ws.Range("N134").Value = -853527.6000000001
$ws.Range("H136").Value = 3303.6875
$ws.Range("I136").Value = 1965
$ws.Range("J136").Value = 5534.8335
$ws.Range("K136").Value = 5895
$ws.Range("L136").Value = 16604.5005
$ws.Range("M136").Value = -3345
$ws.Range("N136").Value = -21704.5005

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5031
$ws.Range("J5").Value = 1463.3334
$ws.Range("L5").Value = 4390.0002
$ws.Range("N5").Value = -4614.0002
$ws.Range("H12").Value = 3788054.8
$ws.Range("I12").Value = 196.55556
$ws.Range("J12").Value = 6410418
$ws.Range("K12").Value = 589.66668
$ws.Range("L12").Value = 19231254
$ws.Range("M12").Value = -416.66668
$ws.Range("N12").Value = -19231600
$ws.Range("H81").Value = 2120.5715
$ws.Range("J81").Value = 3011
$ws.Range("L81").Value = 9033
$ws.Range("N81").Value = -11279
$ws.Range("H84").Value = 2120.5715
$ws.Range("J84").Value = 3011
$ws.Range("L84").Value = 27099
$ws.Range("N84").Value = -38331
$ws.Range("H121").Value = 3177.9268
$ws.Range("I121").Value = 278.91666
$ws.Range("J121").Value = 4377.517
$ws.Range("K121").Value = 836.7499799999999
$ws.Range("L121").Value = 13132.551
$ws.Range("M121").Value = 473.2500200000001
$ws.Range("N121").Value = -15752.551
$ws.Range("H122").Value = 4860.04
$ws.Range("I122").Value = 187.16667
$ws.Range("J122").Value = 9173.462
$ws.Range("K122").Value = 1684.50003
$ws.Range("L122").Value = 82561.158
$ws.Range("M122").Value = 765.4999699999998
$ws.Range("N122").Value = -87461.158
$ws.Range("H131").Value = 841.5599999999999
$ws.Range("J131").Value = 859.1158
$ws.Range("L131").Value = 2577.3474
$ws.Range("N131").Value = -12657.3474
$ws.Range("H135").Value = 5031
$ws.Range("J135").Value = 1463.3334
$ws.Range("L135").Value = 13170.0006
$ws.Range("N135").Value = -18240.0006

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 90912744
$ws.Range("I132").Value = 250001700
$ws.Range("K132").Value = 750005100
$ws.Range("M132").Value = -750002570

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 85836.086
$ws.Range("I122").Value = 113514.78
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 340544.34
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -338094.34
$ws.Range("N122").Value = -13300
$ws.Range("H132").Value = 2321.9473
$ws.Range("I132").Value = 2077.3794
$ws.Range("J132").Value = 3110
$ws.Range("K132").Value = 6232.138199999999
$ws.Range("L132").Value = 9330
$ws.Range("M132").Value = -3702.138199999999
$ws.Range("N132").Value = -14390
$ws.Range("H136").Value = 2006.138
$ws.Range("I136").Value = 1736.2273
$ws.Range("J136").Value = 2854.4285
$ws.Range("K136").Value = 5208.6819
$ws.Range("L136").Value = 8563.2855
$ws.Range("M136").Value = -2658.6819
$ws.Range("N136").Value = -13663.2855

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10346
$ws.Range("H126").Value = 2102544.5
$ws.Range("I126").Value = 3269044.2
$ws.Range("J126").Value = 2844.6
$ws.Range("K126").Value = 9807132.600000001
$ws.Range("L126").Value = 8533.799999999999
$ws.Range("M126").Value = -9804662.600000001
$ws.Range("N126").Value = -13473.8
$ws.Range("H132").Value = 1978011.5
$ws.Range("I132").Value = 2289724
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 6869172
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -6866642
$ws.Range("N132").Value = -16559
